$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.379.50"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.015.74"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.89"
$ws.Range("E5").Value = "  +5.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.41"
$ws.Range("E8").Value = "  -6.18%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.31"
$ws.Range("E12").Value = "  -5.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.312.80"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.805"
$ws.Range("E14").Value = "  -4.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.03"
$ws.Range("E15").Value = "  -6.49%  "
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.018.12"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.330.72"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.93"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0839"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.14"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.65"
$ws.Range("E23").Value = "  +6.65%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.71"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("E27").Value = "  -4.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0649"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  +4.14%  "
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0936"
$ws.Range("E43").Value = "  -5.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.407.07"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.43"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.82"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.204.95"
$ws.Range("E50").Value = "  +0.60%  "
$ws.Range("E51").Value = "  -6.63%  "
